# Update "want to go" counts (column F) for a handful of rows.
# The same update is applied identically to both the "展览" sheet
# and the "全部类型" sheet (they carry duplicate data).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F17").Value = 57
    $ws.Range("F21").Value = 1427
    $ws.Range("F38").Value = 3813
    $ws.Range("F39").Value = 3
    $ws.Range("F43").Value = 62
}
